# Update the "想去人数" (want-to-go count) figures in both the "展览"
# and "全部类型" worksheets, which hold duplicate data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 18
    11 = 4234
    18 = 2975
    31 = 393
    33 = 242
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
